# BlackList.xlsx — fix rule-table DRL snippets and refresh the view state.
#
# Content fix: the "Insured" variable binding/usage in the decision-table
# condition/action columns was written with an inconsistent capitalisation
# ($Insured vs $insured) and the action was missing its trailing semicolon.
# Normalise both to the lower-case binding used elsewhere and terminate the
# action statement.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "`$insured: Insured"
$ws.Range("C8").Value = "`$insured.setStatus(`$param);"

# View refresh: zoom to 150% and move the selection to C9 (as last left by
# the author before saving).
$excel.ActiveWindow.Zoom = 150
$ws.Range("C9").Select()
